# Adds the STUDENTS and EX-STUDENT worksheets, inserts the data_source /
# company_name columns into PROSPECT (pushing status / remark right), and
# appends the extra1 / extra2 / extra3 columns used on all three sheets.

function Set-HeaderCell($ws, $row, $col, $text) {
    # Re-use the existing bold/bordered header style (style index 1, as
    # already used by PROSPECT!A1) by copying that cell's formatting onto
    # the destination, then overwrite the copied text. $script:styleSource
    # is always PROSPECT!A1 regardless of which sheet is being written to.
    $script:styleSource.Copy($ws.Cells.Item($row, $col))
    $ws.Cells.Item($row, $col).Value = $text
}

function Set-TextCell($ws, $row, $col, $text) {
    # Force genuinely text values (including all-digit strings such as
    # phone numbers) to be stored as text rather than being re-interpreted
    # as numbers.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$text
}

function Set-NumCell($ws, $row, $col, $number) {
    $ws.Cells.Item($row, $col).Value = $number
}

$wb = $excel.ActiveWorkbook
$prospect = $wb.Worksheets.Item("PROSPECT")

# Captured once, before any column insert reshuffles PROSPECT - A1 always
# carries the bold/centered/bordered header style (style index 1).
$script:styleSource = $prospect.Range("A1")

# ---------------------------------------------------------------------
# 1. Add the two new worksheets, in order, after PROSPECT.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$students = $wb.Worksheets.Add($null, $lastSheet)
$students.Name = "STUDENTS"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$exStudent = $wb.Worksheets.Add($null, $lastSheet)
$exStudent.Name = "EX-STUDENT"

# ---------------------------------------------------------------------
# 2. PROSPECT: insert "data_source" / "company_name" columns right after
#    "sector" (this shifts the existing "status" / "remark" columns from
#    G/H to I/J), then append extra1 / extra2 / extra3 at the end.
# ---------------------------------------------------------------------
$prospect.Range("G1:H1").EntireColumn.Insert()
Set-HeaderCell $prospect 1 7 "data_source"
Set-HeaderCell $prospect 1 8 "company_name"

Set-HeaderCell $prospect 1 11 "extra1"
Set-HeaderCell $prospect 1 12 "extra2"
Set-HeaderCell $prospect 1 13 "extra3"

# ---------------------------------------------------------------------
# 3. STUDENTS sheet.
# ---------------------------------------------------------------------
Set-HeaderCell $students 1 1 "id"
Set-HeaderCell $students 1 2 "fullname"
Set-HeaderCell $students 1 3 "email"
Set-HeaderCell $students 1 4 "phone"
Set-HeaderCell $students 1 5 "location"
Set-HeaderCell $students 1 6 "dob"
Set-HeaderCell $students 1 7 "courses"
Set-HeaderCell $students 1 8 "registration_fee"
Set-HeaderCell $students 1 9 "tutorial_fee"
Set-HeaderCell $students 1 10 "course_fee"
Set-HeaderCell $students 1 11 "payment_1"
Set-HeaderCell $students 1 12 "payment_2"
Set-HeaderCell $students 1 13 "payment_3"
Set-HeaderCell $students 1 14 "balance"
Set-HeaderCell $students 1 15 "exam"
Set-HeaderCell $students 1 16 "remark_1"
Set-HeaderCell $students 1 17 "remark_2"
Set-HeaderCell $students 1 18 "extra1"
Set-HeaderCell $students 1 19 "extra2"
Set-HeaderCell $students 1 20 "extra3"

Set-NumCell  $students 2 1 2
Set-NumCell  $students 2 2 "Timothy Ojo"
Set-NumCell  $students 2 3 "joeladewole3@gmail.com"
Set-TextCell $students 2 4 "7059575819"
Set-NumCell  $students 2 5 "Ikeja Lagos Nigeria"
Set-NumCell  $students 2 7 "IPOI, IOK, KNIC"
Set-NumCell  $students 2 8 8537493
Set-NumCell  $students 2 9 343554
Set-NumCell  $students 2 10 35457543
Set-NumCell  $students 2 11 34342
Set-NumCell  $students 2 12 34553
Set-NumCell  $students 2 13 4686546
Set-NumCell  $students 2 14 5558654
Set-NumCell  $students 2 15 "Fair"
Set-NumCell  $students 2 16 "This guy is brilliant"
Set-NumCell  $students 2 17 "He hasnt paid up all the money"

# ---------------------------------------------------------------------
# 4. EX-STUDENT sheet.
# ---------------------------------------------------------------------
Set-HeaderCell $exStudent 1 1 "id"
Set-HeaderCell $exStudent 1 2 "fullname"
Set-HeaderCell $exStudent 1 3 "email"
Set-HeaderCell $exStudent 1 4 "phone"
Set-HeaderCell $exStudent 1 5 "location"
Set-HeaderCell $exStudent 1 6 "courses"
Set-HeaderCell $exStudent 1 7 "balance"
Set-HeaderCell $exStudent 1 8 "results"
Set-HeaderCell $exStudent 1 9 "referral_name"
Set-HeaderCell $exStudent 1 10 "referral_number"
Set-HeaderCell $exStudent 1 11 "referral_email"
Set-HeaderCell $exStudent 1 12 "remark"
Set-HeaderCell $exStudent 1 13 "extra1"
Set-HeaderCell $exStudent 1 14 "extra2"
Set-HeaderCell $exStudent 1 15 "extra3"

Set-NumCell  $exStudent 2 1 1
Set-NumCell  $exStudent 2 2 "Sammy Posh"
Set-NumCell  $exStudent 2 3 "samy@posh.com"
Set-TextCell $exStudent 2 4 "9740878636"
Set-NumCell  $exStudent 2 5 "Ibadan"
Set-NumCell  $exStudent 2 6 "IbSL, Posod"
Set-NumCell  $exStudent 2 7 9900
Set-NumCell  $exStudent 2 8 "Good"
Set-NumCell  $exStudent 2 9 "Joe Boy"
Set-TextCell $exStudent 2 10 "9898857670"
Set-NumCell  $exStudent 2 11 "joe@boy.com"
Set-NumCell  $exStudent 2 12 "Fair"

Set-NumCell  $exStudent 3 1 2
Set-NumCell  $exStudent 3 2 "Timothy Ojo"
Set-NumCell  $exStudent 3 3 "tojo@gmail.com"
Set-TextCell $exStudent 3 4 "7059575819"
Set-NumCell  $exStudent 3 5 "San Francisco ,USA"
Set-NumCell  $exStudent 3 6 "IPOI, IOK, KNIC"
Set-NumCell  $exStudent 3 7 49543
Set-NumCell  $exStudent 3 8 "Fair"
Set-NumCell  $exStudent 3 9 "Timothy Foluso"
Set-TextCell $exStudent 3 10 "9048955737594"
Set-NumCell  $exStudent 3 11 "jorti@kb.com"
Set-NumCell  $exStudent 3 12 "The guy finish well"

$prospect.Select()
